$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$Value
    )
    $range = $Sheet.Range($CellRef)
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        # Value would otherwise be auto-converted to a number by Excel;
        # force it to stay as text, matching the source workbook's inlineStr cells.
        $range.Value = "'" + $Value
    } else {
        $range.Value = $Value
    }
}

Set-CellText $ws "D2" '63.824.07'
Set-CellText $ws "E2" '  -0.08%  '
Set-CellText $ws "D3" '2.734.00'
Set-CellText $ws "E3" '  -0.61%  '
Set-CellText $ws "E4" '  +0.11%  '
Set-CellText $ws "D5" '565.28'
Set-CellText $ws "E5" '  -1.68%  '
Set-CellText $ws "D6" '160.53'
Set-CellText $ws "E6" '  +1.83%  '
Set-CellText $ws "E7" '  +0.01%  '
Set-CellText $ws "E8" '  -0.95%  '
Set-CellText $ws "E9" '  -0.14%  '
Set-CellText $ws "E10" '  +4.10%  '
Set-CellText $ws "D11" '5.62'
Set-CellText $ws "E11" '  -1.69%  '
Set-CellText $ws "D12" '0.376'
Set-CellText $ws "E12" '  -1.49%  '
Set-CellText $ws "D13" '3.219.37'
Set-CellText $ws "E13" '  -0.54%  '
Set-CellText $ws "D14" '26.91'
Set-CellText $ws "E14" '  +1.51%  '
Set-CellText $ws "D15" '63.681.29'
Set-CellText $ws "E15" '  +0.28%  '
Set-CellText $ws "E16" '  -0.84%  '
Set-CellText $ws "D17" '2.743.35'
Set-CellText $ws "E17" '  -0.36%  '
Set-CellText $ws "D18" '12.33'
Set-CellText $ws "E18" '  +1.96%  '
Set-CellText $ws "D19" '4.73'
Set-CellText $ws "E19" '  -1.76%  '
Set-CellText $ws "D20" '355.78'
Set-CellText $ws "E20" '  +0.09%  '
Set-CellText $ws "E21" '  -1.26%  '
Set-CellText $ws "D22" '0.997'
Set-CellText $ws "E22" '  -0.12%  '
Set-CellText $ws "D23" '0.519'
Set-CellText $ws "E23" '  -2.59%  '
Set-CellText $ws "D24" '64.17'
Set-CellText $ws "E24" '  -1.51%  '
Set-CellText $ws "E25" '  -0.30%  '
Set-CellText $ws "D26" '1.00'
Set-CellText $ws "E26" '  +0.10%  '
Set-CellText $ws "D27" '8.36'
Set-CellText $ws "E27" '  -1.11%  '
Set-CellText $ws "E28" '  +0.13%  '
Set-CellText $ws "B29" 'PancakeSwap'
Set-CellText $ws "C29" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText $ws "D29" '1.99'
Set-CellText $ws "E29" '  +2.46%  '
Set-CellText $ws "B30" 'Fetch.AI'
Set-CellText $ws "C30" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-CellText $ws "D30" '1.38'
Set-CellText $ws "E30" '  +11.21%  '
Set-CellText $ws "E31" '  +1.44%  '
Set-CellText $ws "D32" '167.41'
Set-CellText $ws "E32" '  -0.55%  '
Set-CellText $ws "D33" '4.91'
Set-CellText $ws "E33" '  -0.02%  '
Set-CellText $ws "E34" '  +2.23%  '
Set-CellText $ws "D35" '20.04'
Set-CellText $ws "E37" '  +1.06%  '
Set-CellText $ws "D38" '0.978'
Set-CellText $ws "E38" '  -0.90%  '
Set-CellText $ws "D39" '347.46'
Set-CellText $ws "E39" '  +4.96%  '
Set-CellText $ws "D40" '6.29'
Set-CellText $ws "E40" '  +2.30%  '
Set-CellText $ws "E41" '  -1.80%  '
Set-CellText $ws "D42" '38.61'
Set-CellText $ws "E42" '  -0.81%  '
Set-CellText $ws "D43" '21.84'
Set-CellText $ws "E43" '  +1.57%  '
Set-CellText $ws "D44" '0.0582'
Set-CellText $ws "E44" '  -0.72%  '
Set-CellText $ws "D45" '20.88'
Set-CellText $ws "E45" '  -2.72%  '
Set-CellText $ws "E46" '  +0.94%  '
Set-CellText $ws "D47" '0.0251'
Set-CellText $ws "E47" '  -1.05%  '
Set-CellText $ws "B48" 'Stellar'
Set-CellText $ws "C48" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws "D48" '0.0994'
Set-CellText $ws "E48" '  -1.39%  '
Set-CellText $ws "B49" 'Aave'
Set-CellText $ws "C49" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText $ws "D49" '132.46'
Set-CellText $ws "E49" '  -1.63%  '
Set-CellText $ws "E50" '  -0.05%  '
Set-CellText $ws "D51" '11.08'
Set-CellText $ws "E51" '  +0.56%  '
